$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Live prereq testing complete: update the four test-case descriptions to
# call out the specific line numbers that were observed during testing.
$ws.Range("B8").Value = "Input is not accepted.  Multiple instances of courses with no prereq listed on line 2, 3, 6, 13."
$ws.Range("B16").Value = "Input is not accepted.  Number of prerequisites exceeds allowed maximum on line 6."
$ws.Range("B18").Value = "Input is not accepted.  Line 6 contains whitespace error."
$ws.Range("B19").Value = "Input is not accepted.  Duplicate prereq on line 7."

# Leave the sheet scrolled/selected where the author left off reviewing.
$ws.Activate()
[void]$ws.Range("B20").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
